$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = 1
$ws.Range("E5").Value = 2
$ws.Range("E13").Value = 3
$ws.Range("E14").Value = 4
$ws.Range("E15").Value = 5

$ws.Range("E16").Select()
